$d = $word.ActiveDocument

# 1) The "Excitedly..." paragraph used to be split across three runs
#    (the middle one wrapped in proofErr spell-check markers around
#    "Ukpo"). Replace the whole sentence, spanning all three runs, with
#    a single wildcard Find/Replace so Word collapses it back down to
#    one run and drops the now-unneeded w:proofErr elements.
$old = "Excitedly, in my current role as President General of Ukpo Improvement Union, I am gaining new experience on fundraising, security management and human capital management."
$pattern = "Excitedly, in my current role as President General of*management."
$d.Content.Find.Execute($pattern, $true, $false, $true, $false, $false, $true, 1, $false, $old, 2)

# 2) After that paragraph, insert a blank paragraph followed by a new
#    paragraph with the "Testify" sentence.
$rng = $d.Paragraphs(4).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.InsertParagraphAfter()

$d.Paragraphs(6).Range.InsertAfter("Testify is also changing my status…")
